# 3 OCTOBER ATTENDANCE - fill in attendance for the 2nd class date (column E),
# mirroring what was entered for the 1st class date (column D), and fix a
# mis-entered value in D18.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 11: header/time row (text, not a numeric attendance mark) ---
$ws.Range("E11").Value = "3:00 PM To 06:00 PM"

# --- Row 12: "Total Classes" row ---
$ws.Range("E12").Value = 3

# --- Correct a mis-entered first-date mark before mirroring it ---
$ws.Range("D18").Value = 3

# --- Student attendance rows (14-42), column E mirrors column D ... ---
$ws.Range("E14").Value = 3
$ws.Range("E15").Value = 0
$ws.Range("E16").Value = 0
$ws.Range("E17").Value = 3
$ws.Range("E18").Value = 3
$ws.Range("E19").Value = 3
$ws.Range("E20").Value = 0
$ws.Range("E21").Value = 3
$ws.Range("E22").Value = 3
$ws.Range("E23").Value = 3
$ws.Range("E24").Value = 3
$ws.Range("E25").Value = 3
$ws.Range("E26").Value = 3
$ws.Range("E27").Value = 3
$ws.Range("E28").Value = 3
$ws.Range("E29").Value = 0
$ws.Range("E30").Value = 3
$ws.Range("E31").Value = 3
$ws.Range("E32").Value = 3
$ws.Range("E33").Value = 3
$ws.Range("E34").Value = 3
$ws.Range("E35").Value = 3
$ws.Range("E36").Value = 3
# ... except row 37, where the student was absent for the 2nd date
$ws.Range("E37").Value = 0
$ws.Range("E38").Value = 3
$ws.Range("E39").Value = 0
$ws.Range("E40").Value = 0
$ws.Range("E41").Value = 0
$ws.Range("E42").Value = 0

# --- Leave the view scrolled/selected where the user last worked ---
$ws.Range("D27").Select()

Write-Output "attendance for 2nd class date filled in"
